$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price (D) column so that values such as
# "1.000" or "0.4517" are stored as literal text instead of being
# auto-converted to numbers by Excel (matches original inline-string data).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.324.34'
$ws.Range('E2').Value = '  -4.53%  '
$ws.Range('D3').Value = '1.858.62'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').Value = '322.43'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').Value = '0.4517'
$ws.Range('E7').Value = '  -6.51%  '
$ws.Range('D8').Value = '0.3857'
$ws.Range('E8').Value = '  -5.43%  '
$ws.Range('D9').Value = '48.18'
$ws.Range('E9').Value = '  -11.06%  '
$ws.Range('D10').Value = '0.07889'
$ws.Range('E10').Value = '  -7.42%  '
$ws.Range('D11').Value = '1.021'
$ws.Range('E11').Value = '  -4.16%  '
$ws.Range('D12').Value = '21.38'
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('D13').Value = '1.861.35'
$ws.Range('E13').Value = '  -7.06%  '
$ws.Range('D14').Value = '7.159'
$ws.Range('E14').Value = '  -6.25%  '
$ws.Range('D15').Value = '5.874'
$ws.Range('E15').Value = '  -5.44%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('E17').Value = '  -4.23%  '
$ws.Range('D18').Value = '85.44'
$ws.Range('E18').Value = '  -6.30%  '
$ws.Range('D19').Value = '0.06524'
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('D20').Value = '17.05'
$ws.Range('E20').Value = '  -8.57%  '
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').Value = '5.513'
$ws.Range('E22').Value = '  -6.50%  '
$ws.Range('D23').Value = '27.328.18'
$ws.Range('E23').Value = '  -4.66%  '
$ws.Range('E24').Value = '  -6.79%  '
$ws.Range('D25').Value = '2.269'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = '2.081.41'
$ws.Range('E26').Value = '  -7.01%  '
$ws.Range('D27').Value = '151.82'
$ws.Range('E27').Value = '  -3.14%  '
$ws.Range('D28').Value = '19.72'
$ws.Range('E28').Value = '  -3.41%  '
$ws.Range('D29').Value = '2.061'
$ws.Range('E29').Value = '  -6.06%  '
$ws.Range('D30').Value = '5.504'
$ws.Range('E30').Value = '  -7.17%  '
$ws.Range('D31').Value = '120.39'
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('D32').Value = '1.476'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = '0.09298'
$ws.Range('E33').Value = '  -4.29%  '
$ws.Range('D34').Value = '0.9372'
$ws.Range('E34').Value = '  -6.03%  '
$ws.Range('D35').Value = '3.596'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').Value = '5.276'
$ws.Range('E36').Value = '  -6.83%  '
$ws.Range('D37').Value = '0.02228'
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('D38').Value = '0.05988'
$ws.Range('E38').Value = '  -4.51%  '
$ws.Range('D39').Value = '1.213'
$ws.Range('E39').Value = '  -3.68%  '
$ws.Range('D40').Value = '8.270'
$ws.Range('E40').Value = '  -10.05%  '
$ws.Range('D41').Value = '0.9995'
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('D42').Value = '0.5911'
$ws.Range('D43').Value = '0.1889'
$ws.Range('E43').Value = '  -1.87%  '
$ws.Range('D44').Value = '10.13'
$ws.Range('E44').Value = '  -10.30%  '
$ws.Range('D45').Value = '1.262'
$ws.Range('E45').Value = '  -6.86%  '
$ws.Range('D46').Value = '0.5629'
$ws.Range('E46').Value = '  -5.98%  '
$ws.Range('D47').Value = '12.00'
$ws.Range('E47').Value = '  -8.44%  '
$ws.Range('D48').Value = '1.925'
$ws.Range('E48').Value = '  -7.49%  '
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('D50').Value = '0.06803'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').Value = '108.26'
$ws.Range('E51').Value = '  -3.24%  '
